$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'65.903.98"
$ws.Range("E2").Value = "  -2.12%  "

# Row 3
$ws.Range("D3").Value = "'3.411.69"
$ws.Range("E3").Value = "  -1.54%  "

# Row 4
$ws.Range("E4").Value = "  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'582.05"
$ws.Range("E5").Value = "  -1.91%  "

# Row 6
$ws.Range("D6").Value = "'171.93"
$ws.Range("E6").Value = "  -4.69%  "

# Row 8
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  -3.92%  "

# Row 9
$ws.Range("D9").Value = "'3.413.72"
$ws.Range("E9").Value = "  -1.39%  "

# Row 10
$ws.Range("D10").Value = "'0.129"
$ws.Range("E10").Value = "  -7.35%  "

# Row 11
$ws.Range("D11").Value = "'6.82"
$ws.Range("E11").Value = "  -1.91%  "

# Row 12
$ws.Range("D12").Value = "'0.407"
$ws.Range("E12").Value = "  -5.24%  "

# Row 13
$ws.Range("D13").Value = "'4.004.81"
$ws.Range("E13").Value = "  -1.33%  "

# Row 14
$ws.Range("E14").Value = "  -0.71%  "

# Row 15
$ws.Range("D15").Value = "'29.72"
$ws.Range("E15").Value = "  -7.47%  "

# Row 16
$ws.Range("D16").Value = "'65.999.59"
$ws.Range("E16").Value = "  -1.95%  "

# Row 17
$ws.Range("D17").Value = "'0.0000169"
$ws.Range("E17").Value = "  -4.36%  "

# Row 18
$ws.Range("D18").Value = "'3.419.32"
$ws.Range("E18").Value = "  -1.18%  "

# Row 19
$ws.Range("D19").Value = "'5.86"
$ws.Range("E19").Value = "  -5.55%  "

# Row 20
$ws.Range("D20").Value = "'13.63"
$ws.Range("E20").Value = "  -3.50%  "

# Row 21
$ws.Range("D21").Value = "'364.59"
$ws.Range("E21").Value = "  -7.44%  "

# Row 22
$ws.Range("D22").Value = "'7.64"
$ws.Range("E22").Value = "  -3.56%  "

# Row 23
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  +0.05%  "

# Row 24
$ws.Range("E24").Value = "  -1.45%  "

# Row 25
$ws.Range("D25").Value = "'70.98"
$ws.Range("E25").Value = "  -1.01%  "

# Row 26
$ws.Range("D26").Value = "'0.523"
$ws.Range("E26").Value = "  -3.01%  "

# Row 27
$ws.Range("E27").Value = "  -3.09%  "

# Row 28
$ws.Range("D28").Value = "'9.56"
$ws.Range("E28").Value = "  -7.69%  "

# Row 29
$ws.Range("E29").Value = "  +0.73%  "

# Row 30
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.06%  "

# Row 31
$ws.Range("D31").Value = "'5.75"
$ws.Range("E31").Value = "  -6.05%  "

# Row 32
$ws.Range("D32").Value = "'23.50"
$ws.Range("E32").Value = "  -0.09%  "

# Row 33
$ws.Range("D33").Value = "'1.97"
$ws.Range("E33").Value = "  -3.79%  "

# Row 34
$ws.Range("E34").Value = "  -0.06%  "

# Row 35
$ws.Range("D35").Value = "'1.28"
$ws.Range("E35").Value = "  -8.77%  "

# Row 36
$ws.Range("D36").Value = "'7.00"
$ws.Range("E36").Value = "  -4.46%  "

# Row 37
$ws.Range("E37").Value = "  -3.57%  "

# Row 38
$ws.Range("D38").Value = "'160.60"
$ws.Range("E38").Value = "  -0.36%  "

# Row 39
$ws.Range("D39").Value = "'28.81"
$ws.Range("E39").Value = "  +10.23%  "

# Row 40
$ws.Range("D40").Value = "'0.876"
$ws.Range("E40").Value = "  -0.33%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.75"
$ws.Range("E41").Value = "  -6.78%  "

# Row 42
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D42").Value = "'2.57"
$ws.Range("E42").Value = "  -9.23%  "

# Row 43
$ws.Range("D43").Value = "'2.695.28"
$ws.Range("E43").Value = "  -2.26%  "

# Row 44
$ws.Range("D44").Value = "'4.37"
$ws.Range("E44").Value = "  -6.32%  "

# Row 45
$ws.Range("E45").Value = "  -6.16%  "

# Row 46
$ws.Range("D46").Value = "'0.0676"
$ws.Range("E46").Value = "  -6.04%  "

# Row 47
$ws.Range("D47").Value = "'39.86"
$ws.Range("E47").Value = "  -3.67%  "

# Row 48
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'23.92"
$ws.Range("E48").Value = "  -9.08%  "

# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0287"
$ws.Range("E49").Value = "  -3.69%  "

# Row 50
$ws.Range("D50").Value = "'304.24"
$ws.Range("E50").Value = "  -6.64%  "

# Row 51
$ws.Range("D51").Value = "'0.811"
$ws.Range("E51").Value = "  -3.61%  "
